{"js": "// Fix final encoding errors: stray literal \">=\" sequences left in front of\n// mis-encoded special characters (superscript numerals, em dash, bullet/\n// approx signs) by a previous cleanup pass. Each fix below targets the\n// full text of one existing run so the run/rPr structure is left untouched.\nconst fixes = [\n  [\"Dwi Anggriani>=\u00b9, Syaiful Bachri Mustamin>=\u00b9\", \"Dwi Anggriani\u00b9, Syaiful Bachri Mustamin\u00b9\"],\n  [\", Muhammad Atnang>=\u00b9, Kartini Aprilia Pratiwi Nuzry>=\u00b9\", \", Muhammad Atnang\u00b9, Kartini Aprilia Pratiwi Nuzry\u00b9\"],\n  [\">=\u00b9Department of Information Technology, Institut Sains Teknologi dan Kesehatan \u2019Aisyiyah Kendari, Kendari, Indonesia\", \"\u00b9Department of Information Technology, Institut Sains Teknologi dan Kesehatan \u2019Aisyiyah Kendari, Kendari, Indonesia\"],\n  [\": Using majority voting (>=\u2030\u00a52/4 runs), we calculated accuracy, sensitivity, specificity, precision, and F1-score\", \": Using majority voting (\u2030\u00a52/4 runs), we calculated accuracy, sensitivity, specificity, precision, and F1-score\"],\n  [\"for nearly all cases, with true negatives >=\u2030\u02c60.\", \"for nearly all cases, with true negatives \u2248 0.\"],\n  [\"High consistency indicates LLMs reliably apply learned reasoning patterns>=\u20ac\u201dthey are systematically biased rather than randomly erring. This\", \"High consistency indicates LLMs reliably apply learned reasoning patterns\u20ac\u201dthey are systematically biased rather than randomly erring. This\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of fixes) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Pattern not found: \" + oldText);\n  }\n\n  // Each broken string is unique in the document, so replace every hit\n  // (normally exactly one) with the corrected text.\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n", "ps1": "# Fix final encoding errors: stray literal \">=\" sequences left in front of\n# mis-encoded special characters (superscript numerals, em dash, bullet/\n# approx signs) by a previous cleanup pass. Each Find/Replace below targets\n# the full text of one existing run, so Word only rewrites that run's text.\n\n$d = $word.ActiveDocument\n\n$fixes = @(\n    ,@(\"Dwi Anggriani>=\u00b9, Syaiful Bachri Mustamin>=\u00b9\", \"Dwi Anggriani\u00b9, Syaiful Bachri Mustamin\u00b9\")\n    ,@(\", Muhammad Atnang>=\u00b9, Kartini Aprilia Pratiwi Nuzry>=\u00b9\", \", Muhammad Atnang\u00b9, Kartini Aprilia Pratiwi Nuzry\u00b9\")\n    ,@(\">=\u00b9Department of Information Technology, Institut Sains Teknologi dan Kesehatan \u2019Aisyiyah Kendari, Kendari, Indonesia\", \"\u00b9Department of Information Technology, Institut Sains Teknologi dan Kesehatan \u2019Aisyiyah Kendari, Kendari, Indonesia\")\n    ,@(\": Using majority voting (>=\u2030\u00a52/4 runs), we calculated accuracy, sensitivity, specificity, precision, and F1-score\", \": Using majority voting (\u2030\u00a52/4 runs), we calculated accuracy, sensitivity, specificity, precision, and F1-score\")\n    ,@(\"for nearly all cases, with true negatives >=\u2030\u02c60.\", \"for nearly all cases, with true negatives \u2248 0.\")\n    ,@(\"High consistency indicates LLMs reliably apply learned reasoning patterns>=\u20ac\u201dthey are systematically biased rather than randomly erring. This\", \"High consistency indicates LLMs reliably apply learned reasoning patterns\u20ac\u201dthey are systematically biased rather than randomly erring. This\")\n)\n\nforeach ($fix in $fixes) {\n    $oldText = $fix[0]\n    $newText = $fix[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"Pattern not found: $oldText\"\n    }\n}\n\n"}
